# FSC_UN_metrics.xlsx edit: insert a "runtime" column between the
# original 3-metric block (B:D) and the standard_blocking 3-metric
# block (old E:G), and append a trailing "runtime" column after the
# standard_blocking block.
#
# Layout before:  A | B C D (original: recall,precision,f1) | E F G (standard_blocking: recall,precision,f1)
# Layout after:   A | B C D E (original: recall,precision,f1,runtime) | F G H I (standard_blocking: recall,precision,f1,runtime)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E; this shifts old E:G -> F:H, and also shifts
# the merged header cells along with them automatically (mergeCell
# refs stay the same width, just slide right).
$ws.Columns("E:E").Insert()

# Grow the two header merges by one column each so they keep spanning
# their (now 4-wide) metric blocks.
$ws.Range("B1:E1").MergeCells = $true
$ws.Range("F1:I1").MergeCells = $true

# Re-merging above resets the border styling on individual cells in
# the merged band (Excel gives edge/middle cells different border
# variants). Stamp the original header style (taken from A1, which was
# never touched by the merge calls) back over the whole header row so
# every header cell uses the same uniform style, like before.
$ws.Range("A1").Copy()
$ws.Range("B1:I1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Row 2 header labels: new "runtime" column headers for each block.
# I2 is a brand-new cell with no style yet -- borrow D2's header style
# (D2 sits outside any merged range, so it was never touched above)
# before writing its value.
$ws.Range("D2").Copy()
$ws.Range("I2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E2").Value = "runtime"
$ws.Range("I2").Value = "runtime"

# Per-model runtime values (mm:ss text) for the new E column (original
# block) and the new I column (standard_blocking block).
$runtimes = @{
  4  = @("38:05", "17:25")
  5  = @("17:52", "07:28")
  6  = @("18:02", "07:38")
  7  = @("15:30", "06:19")
  8  = @("45:44", "38:31")
  9  = @("17:56", "15:20")
  10 = @("18:18", "15:41")
  11 = @("13:10", "11:51")
  12 = @("40:46", "37:31")
  13 = @("16:00", "14:37")
  14 = @("16:27", "15:00")
  15 = @("11:48", "10:53")
  16 = @("46:37", "64:01")
  17 = @("18:20", "23:14")
  18 = @("19:05", "24:03")
  19 = @("14:18", "19:30")
  20 = @("55:33", "48:06")
  21 = @("22:51", "19:54")
  22 = @("23:16", "20:16")
  23 = @("17:01", "15:00")
}

foreach ($r in $runtimes.Keys) {
  $pair = $runtimes[$r]
  $ws.Range("E$r").Value = $pair[0]
  $ws.Range("I$r").Value = $pair[1]
}

Write-Output ("dimension=" + $ws.UsedRange.Address())
